$d = $word.ActiveDocument

# Remove the <w:contextualSpacing w:val="0"/> element from every paragraph's
# paragraph-properties. The Word object model does not expose a
# ParagraphFormat.ContextualSpacing property in this runtime, so we round-trip
# each paragraph through its WordOpenXML representation, strip the element
# from the raw OOXML, and re-insert it with InsertXML (which replaces only
# the content of that specific range).
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $xml = $r.WordOpenXML
    if ($xml -match "<w:contextualSpacing\b[^/]*/>") {
        $newXml = $xml -replace "<w:contextualSpacing\b[^/]*/>", ""
        [void]$r.InsertXML($newXml)
    }
}
